# "Yard does not work on Excel"
# Add three new "Yard" placeholder blocks (Section Y, block numbers 0-2) right
# after the existing track data, and relocate the MIN/seconds-per-block
# summary formulas that used to live at the bottom of column K out to column
# M (since column K is now back in use by the new block rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the trailing summary formulas out of K153:K154 into M157:M158
#        so the new data rows below can reuse K153/K154. ---
$ws.Cells.Item(153, 11).ClearContents()
$ws.Cells.Item(154, 11).ClearContents()

$ws.Cells.Item(157, 13).Formula = "=MIN(K2:K151)"
$ws.Cells.Item(158, 13).Formula = "=M157/1.2"

# --- 2. Append three new "Yard" test rows (152-154), Section Y, with all
#        zeroed-out measurements. ---
$yardRows = 152, 153, 154
$blockNums = 0, 1, 2

for ($i = 0; $i -lt $yardRows.Length; $i++) {
    $r = $yardRows[$i]
    $ws.Cells.Item($r, 1).Value = "Green"
    $ws.Cells.Item($r, 2).Value = "Y"
    $ws.Cells.Item($r, 3).Value = $blockNums[$i]
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Formula = "=E" + $r + "*D" + $r + "/100"
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
}

# --- 3. Match the sheet view seen in the edited workbook: scrolled down to
#        the newly added rows, with the header row still frozen, and the
#        active cell sitting on the new Yard section. ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C155").Select() | Out-Null
